$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2 through 407) from 45171 (2023-09-02) to 45172 (2023-09-03).
for ($i = 2; $i -le 407; $i++) {
    $ws.Cells.Item($i, 3).Value2 = 45172
}

# Rows 406 and 407 also swapped their "Beteckning" (A) and "Area (ha)" (G)
# values - effectively the two records traded places while keeping the
# rest of their row data (B, D, E, H:Q) in place.
$a406 = $ws.Cells.Item(406, 1).Value2
$a407 = $ws.Cells.Item(407, 1).Value2
$g406 = $ws.Cells.Item(406, 7).Value2
$g407 = $ws.Cells.Item(407, 7).Value2

$ws.Cells.Item(406, 1).Value2 = $a407
$ws.Cells.Item(407, 1).Value2 = $a406
$ws.Cells.Item(406, 7).Value2 = $g407
$ws.Cells.Item(407, 7).Value2 = $g406
